$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.5121479034423828
$ws.Range("D2").Value = 67.27800000000001

$ws.Range("C3").Value = 0.1740388870239258
$ws.Range("D3").Value = 38.991

$ws.Range("C4").Value = 0.9380471706390381
$ws.Range("D4").Value = 56.173

$ws.Range("C5").Value = 0.3974390029907227
$ws.Range("D5").Value = 73.245

$ws.Range("C6").Value = 0.1695139408111572
$ws.Range("D6").Value = 6.017

$ws.Range("C7").Value = 0.1692137718200684
$ws.Range("D7").Value = 6.017

$ws.Range("C8").Value = 0.5016591548919678
$ws.Range("D8").Value = 54.218

$ws.Range("C9").Value = 0.4957399368286133
$ws.Range("D9").Value = 71.29000000000001

$ws.Range("C10").Value = 0.1831440925598145
$ws.Range("D10").Value = 7.972

$ws.Range("C11").Value = 0.1768581867218018
$ws.Range("D11").Value = 7.972

$ws.Range("C12").Value = 0.4625730514526367
$ws.Range("D12").Value = 56.173

$ws.Range("C13").Value = 0.404789924621582
$ws.Range("D13").Value = 73.245

$ws.Range("C14").Value = 0.1723949909210205
$ws.Range("D14").Value = 64.398

$ws.Range("C15").Value = 0.1738078594207764
$ws.Range("D15").Value = 522.806

$ws.Range("C16").Value = 0.4329590797424316
$ws.Range("D16").Value = 58.093

$ws.Range("C17").Value = 0.5745940208435059
$ws.Range("D17").Value = 522.668

$ws.Range("C18").Value = 0.1825730800628662
$ws.Range("D18").Value = 6.215

$ws.Range("C19").Value = 0.1899170875549316
$ws.Range("D19").Value = 6.215

$ws.Range("C20").Value = 0.4363090991973877
$ws.Range("D20").Value = 56.138

$ws.Range("C21").Value = 0.5788819789886475
$ws.Range("D21").Value = 520.713

$ws.Range("C22").Value = 0.1883602142333984
$ws.Range("D22").Value = 8.17

$ws.Range("C23").Value = 0.18572998046875
$ws.Range("D23").Value = 8.17

$ws.Range("C24").Value = 0.4326059818267822
$ws.Range("D24").Value = 58.093

$ws.Range("C25").Value = 0.7214400768280029
$ws.Range("D25").Value = 522.668
